# Atualização de bases das ligas, do dia: 14-04-2024 às 18:28
# Swap betting-odds rows whose ids were recorded out of order, and drop a
# stray duplicate row (id 8089991 / B287) that shouldn't be in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($rowA, $rowB) {
    $rangeA = $ws.Range("B$rowA`:AC$rowA")
    $rangeB = $ws.Range("B$rowB`:AC$rowB")
    $valsA = $rangeA.Value2
    $valsB = $rangeB.Value2
    $rangeA.Value2 = $valsB
    $rangeB.Value2 = $valsA
}

# Simple pairwise swaps (columns B..AC only; column A "index" stays put).
Swap-Rows 47 48
Swap-Rows 55 56
Swap-Rows 71 72
Swap-Rows 133 134
Swap-Rows 213 214
Swap-Rows 221 222
Swap-Rows 229 231
Swap-Rows 232 233
Swap-Rows 248 249
Swap-Rows 251 252
Swap-Rows 271 272

# Three-way cyclic rotation: new263 = old264, new264 = old265, new265 = old263
$r263 = $ws.Range("B263:AC263")
$r264 = $ws.Range("B264:AC264")
$r265 = $ws.Range("B265:AC265")
$v263 = $r263.Value2
$v264 = $r264.Value2
$v265 = $r265.Value2
$r263.Value2 = $v264
$r264.Value2 = $v265
$r265.Value2 = $v263

# Drop the last data row (was row 287, id 8089991) entirely.
$ws.Rows(287).Delete()
